# The sheet tracks one row per (market, quality) price observation, ordered
# by date. This commit adds a newer "Región de Ñuble" observation pair
# (Primera/Segunda) for 2022-05-24 (serial 44705) at the top of the
# Acelga / Vega Monumental Concepción block, pushing the existing rows
# 128.. down by two (dimension grows from R259 to R261).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 128; this shifts the
# existing data rows 128..259 down to 130..261 automatically.
$ws.Rows.Item(128).Insert()
$ws.Rows.Item(128).Insert()

# New row 128 ("Primera" observation for the newly added date)
$ws.Range("A128").Value = 11
$ws.Range("B128").Value = "Vega Monumental Concepción"
$ws.Range("C128").Value = "Bíobío"
$ws.Range("D128").Value = 44705
$ws.Range("E128").Value = 8
$ws.Range("F128").Value = 100112009
$ws.Range("G128").Value = "Acelga"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 200
$ws.Range("K128").Value = 600
$ws.Range("L128").Value = 700
$ws.Range("M128").Value = 650
$ws.Range("N128").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O128").Value = "Región de Ñuble"
$ws.Range("P128").Value = 650
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"

# New row 129 ("Segunda" observation for the newly added date)
$ws.Range("A129").Value = 11
$ws.Range("B129").Value = "Vega Monumental Concepción"
$ws.Range("C129").Value = "Bíobío"
$ws.Range("D129").Value = 44705
$ws.Range("E129").Value = 8
$ws.Range("F129").Value = 100112009
$ws.Range("G129").Value = "Acelga"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 100
$ws.Range("K129").Value = 500
$ws.Range("L129").Value = 500
$ws.Range("M129").Value = 500
$ws.Range("N129").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O129").Value = "Región de Ñuble"
$ws.Range("P129").Value = 500
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"
